$wb = $excel.ActiveWorkbook

# --- Update "scraped_at" timestamps (column K) on the snapshot sheet ---
$snapshot = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    "K2" = "2025-11-17T11:11:06.343361+00:00"
    "K3" = "2025-11-17T11:11:06.343381+00:00"
    "K4" = "2025-11-17T11:11:06.343390+00:00"
    "K5" = "2025-11-17T11:11:08.296613+00:00"
    "K6" = "2025-11-17T11:11:08.296631+00:00"
    "K7" = "2025-11-17T11:11:08.296640+00:00"
    "K8" = "2025-11-17T11:11:10.950508+00:00"
    "K9" = "2025-11-17T11:11:13.294859+00:00"
    "K10" = "2025-11-17T11:11:15.987370+00:00"
    "K11" = "2025-11-17T11:11:15.987401+00:00"
    "K12" = "2025-11-17T11:11:21.029493+00:00"
    "K13" = "2025-11-17T11:11:23.422784+00:00"
    "K14" = "2025-11-17T11:11:26.100105+00:00"
    "K15" = "2025-11-17T11:11:26.100133+00:00"
    "K16" = "2025-11-17T11:11:26.100151+00:00"
    "K17" = "2025-11-17T11:11:28.381254+00:00"
    "K18" = "2025-11-17T11:11:30.754090+00:00"
    "K19" = "2025-11-17T11:11:30.754178+00:00"
    "K20" = "2025-11-17T11:11:33.057939+00:00"
    "K21" = "2025-11-17T11:11:35.385590+00:00"
    "K22" = "2025-11-17T11:11:35.385624+00:00"
    "K23" = "2025-11-17T11:11:35.385645+00:00"
    "K24" = "2025-11-17T11:11:35.385663+00:00"
    "K25" = "2025-11-17T11:11:37.743881+00:00"
    "K26" = "2025-11-17T11:11:37.743912+00:00"
    "K27" = "2025-11-17T11:11:40.043951+00:00"
    "K28" = "2025-11-17T11:11:40.043980+00:00"
    "K29" = "2025-11-17T11:11:40.044000+00:00"
    "K30" = "2025-11-17T11:11:42.401593+00:00"
    "K31" = "2025-11-17T11:11:42.401624+00:00"
    "K32" = "2025-11-17T11:11:44.804044+00:00"
    "K33" = "2025-11-17T11:11:44.804074+00:00"
    "K34" = "2025-11-17T11:11:44.804092+00:00"
    "K35" = "2025-11-17T11:11:44.804109+00:00"
    "K36" = "2025-11-17T11:11:44.804128+00:00"
    "K37" = "2025-11-17T11:11:47.168150+00:00"
    "K38" = "2025-11-17T11:11:47.168180+00:00"
    "K39" = "2025-11-17T11:11:51.924215+00:00"
    "K40" = "2025-11-17T11:11:51.924245+00:00"
    "K41" = "2025-11-17T11:11:51.924265+00:00"
    "K42" = "2025-11-17T11:11:51.924283+00:00"
    "K43" = "2025-11-17T11:11:54.214357+00:00"
}

foreach ($cellRef in $timestamps.Keys) {
    $snapshot.Range($cellRef).Value = $timestamps[$cellRef]
}

# --- Remove the processed new-injury row from the new_injured sheet ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()

